# "Add files via upload" re-save: the sheets (previously named after each
# student) are generalized to data1..data10, and the view is left parked on
# the last sheet (which becomes the new active tab / selection), while the
# sheet that used to be active (Shin/data7) keeps its last selection but is
# no longer the active tab.

$wb = $excel.ActiveWorkbook

$newNames = @("data1","data2","data3","data4","data5","data6","data7","data8","data9","data10")

for ($i = 0; $i -lt $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i + 1)
    $ws.Name = $newNames[$i]
}

# Sheet 7 ("data7", formerly "Shin") was the active tab with K22 selected;
# revisit it to leave its own selection as-is, then move off of it so it's
# no longer the active tab.
$ws7 = $wb.Worksheets.Item(7)
$ws7.Activate()
$ws7.Range("K22").Select()

# Sheet 10 ("data10", formerly "Laughlin") becomes the newly active tab,
# with B6 selected.
$ws10 = $wb.Worksheets.Item(10)
$ws10.Activate()
$ws10.Range("B6").Select()
